$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.6
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 2.82
$ws.Range("L4").Value = 3.4
$ws.Range("N4").Value = 6.9
$ws.Range("O4").Value = 1.32
$ws.Range("P4").Value = 3.1
$ws.Range("Q4").Value = 1.98
$ws.Range("R4").Value = 1.78
$ws.Range("S4").Value = 1.39
$ws.Range("T4").Value = 2.77
$ws.Range("U4").Value = 1.72
$ws.Range("V4").Value = 2
$ws.Range("Z4").Value = 25
$ws.Range("AC4").Value = 6.9
$ws.Range("AD4").Value = 6.1
$ws.Range("AE4").Value = 13
$ws.Range("AG4").Value = 9
$ws.Range("AM4").Value = 450
$ws.Range("AT4").Value = 2.77
$ws.Range("AW4").Value = 4.85
$ws.Range("AX4").Value = 15.5
$ws.Range("AZ4").Value = 70
$ws.Range("H5").Value = 4.6
$ws.Range("I5").Value = 1.31
$ws.Range("J5").Value = 7.5
$ws.Range("K5").Value = 2.42
$ws.Range("N5").Value = 8.5
$ws.Range("P5").Value = 3.9
$ws.Range("Q5").Value = 1.65
$ws.Range("R5").Value = 2.15
$ws.Range("S5").Value = 1.33
$ws.Range("T5").Value = 3.05
$ws.Range("U5").Value = 2.02
$ws.Range("W5").Value = 24
$ws.Range("Z5").Value = 250
$ws.Range("AA5").Value = 110
$ws.Range("AB5").Value = 90
$ws.Range("AC5").Value = 8.5
$ws.Range("AK5").Value = 11.25
$ws.Range("AO5").Value = 50
$ws.Range("AQ5").Value = 400
$ws.Range("AR5").Value = 350
$ws.Range("AT5").Value = 3.05
$ws.Range("AV5").Value = 80
$ws.Range("AX5").Value = 5.8
$ws.Range("AY5").Value = 16.5
$ws.Range("AZ5").Value = 15.5
$ws.Range("I6").Value = 7.7
$ws.Range("L6").Value = 6.8
$ws.Range("O6").Value = 1.16
$ws.Range("P6").Value = 4.65
$ws.Range("Q6").Value = 1.5
$ws.Range("R6").Value = 2.42
$ws.Range("S6").Value = 1.27
$ws.Range("U6").Value = 1.85
$ws.Range("V6").Value = 1.85
$ws.Range("W6").Value = 8.75
$ws.Range("X6").Value = 7.1
$ws.Range("Z6").Value = 8.25
$ws.Range("AA6").Value = 10.5
$ws.Range("AG6").Value = 25
$ws.Range("AI6").Value = 25
$ws.Range("AK6").Value = 80
$ws.Range("AM6").Value = 600
$ws.Range("AR6").Value = 37
$ws.Range("AU6").Value = 8.5
$ws.Range("AW6").Value = 9
$ws.Range("BB6").Value = 500
$ws.Range("G7").Value = 4.8
$ws.Range("H7").Value = 3.5
$ws.Range("I7").Value = 1.65
$ws.Range("J7").Value = 5.1
$ws.Range("K7").Value = 2.18
$ws.Range("L7").Value = 2.2
$ws.Range("N7").Value = 7.1
$ws.Range("O7").Value = 1.32
$ws.Range("P7").Value = 3.1
$ws.Range("Q7").Value = 1.95
$ws.Range("U7").Value = 1.91
$ws.Range("W7").Value = 12.5
$ws.Range("X7").Value = 28
$ws.Range("Y7").Value = 16
$ws.Range("Z7").Value = 90
$ws.Range("AA7").Value = 45
$ws.Range("AC7").Value = 7.1
$ws.Range("AD7").Value = 7
$ws.Range("AG7").Value = 6.3
$ws.Range("AH7").Value = 7.4
$ws.Range("AJ7").Value = 12.5
$ws.Range("AN7").Value = 6.6
$ws.Range("AP7").Value = 35
$ws.Range("AV7").Value = 75
$ws.Range("AW7").Value = 3.45
$ws.Range("AX7").Value = 8
$ws.Range("AY7").Value = 18
$ws.Range("G8").Value = 2.35
$ws.Range("I8").Value = 2.67
$ws.Range("J8").Value = 2.95
$ws.Range("K8").Value = 2.2
$ws.Range("L8").Value = 3.2
$ws.Range("S8").Value = 1.34
$ws.Range("T8").Value = 3
$ws.Range("V8").Value = 2.25
$ws.Range("W8").Value = 9.5
$ws.Range("X8").Value = 13
$ws.Range("Y8").Value = 9.25
$ws.Range("Z8").Value = 25
$ws.Range("AA8").Value = 18
$ws.Range("AB8").Value = 24
$ws.Range("AD8").Value = 6.8
$ws.Range("AE8").Value = 12
$ws.Range("AH8").Value = 16
$ws.Range("AI8").Value = 10
$ws.Range("AJ8").Value = 32
$ws.Range("AK8").Value = 20
$ws.Range("AL8").Value = 24
$ws.Range("AN8").Value = 4.5
$ws.Range("AO8").Value = 12.5
$ws.Range("AP8").Value = 18.5
$ws.Range("AQ8").Value = 50
$ws.Range("AR8").Value = 75
$ws.Range("AS8").Value = 200
$ws.Range("AT8").Value = 3
$ws.Range("AW8").Value = 4.85
$ws.Range("AX8").Value = 14
$ws.Range("AY8").Value = 19
$ws.Range("AZ8").Value = 60
$ws.Range("BA8").Value = 80
$ws.Range("BB8").Value = 200
